$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 156, pushing existing rows 156..197 down to 157..198
$ws.Rows.Item(156).Insert()

# Populate the new row 156 with data (a new Uva price record for
# Agricola del Norte S.A. de Arica / Arica y Parinacota)
$ws.Cells.Item(156, 1).Value = 1
$ws.Cells.Item(156, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(156, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(156, 4).Value = 45275
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(157, 4).NumberFormat
$ws.Cells.Item(156, 5).Value = 15
$ws.Cells.Item(156, 6).Value = "Fruta"
$ws.Cells.Item(156, 7).Value = 100109
$ws.Cells.Item(156, 8).Value = "Uva"
$ws.Cells.Item(156, 9).Value = 100109001
$ws.Cells.Item(156, 10).Value = "Uva"
$ws.Cells.Item(156, 11).Value = "Red Globe"
$ws.Cells.Item(156, 12).Value = "Segunda"
$ws.Cells.Item(156, 13).Value = 250
$ws.Cells.Item(156, 14).Value = 24000
$ws.Cells.Item(156, 15).Value = 25000
$ws.Cells.Item(156, 16).Value = 24500
$ws.Cells.Item(156, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(156, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(156, 19).Value = 2450
$ws.Cells.Item(156, 20).Value = 10
